# Generate Report for Handoff
# b.md has now been handed off again (new handoff files / timestamps),
# so its status flips from "Handed back: in sync with en-US" to
# "Ready for handoff" on every sheet, and the corresponding handoff
# file name + handoff datetime are refreshed on the per-locale sheets.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet: row for b.md (row 3) ----
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"
$ov.Range("D3").Value = "2016-35-13 16:35:11"

# ---- zh-cn sheet: row for b.md (row 3) ----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zh.Range("E3").Value = "2016-03-13 16:35:07"
$zh.Range("D3").Hyperlinks.Item(1).TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"

# ---- de-de sheet: row for b.md (row 3) ----
$de = $wb.Worksheets.Item("de-de")
$de.Range("C3").Value = "Ready for handoff"
$de.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$de.Range("E3").Value = "2016-03-13 16:35:11"
$de.Range("D3").Hyperlinks.Item(1).TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
